# The workbook has a "slider" (Forms scroll-bar control named "Scroll Bar 1")
# on Sheet1 whose linked cell is D3. D2 holds the formula "=D3/10".
# Moving the slider to 34 sets D3 to 34, which ripples into D2 (=3.4).
#
# Drive this the same way a user dragging the scrollbar would: write the
# linked cell's value directly (Excel's automatic-calc engine then updates
# the dependent formula cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set the scroll bar's linked cell to the new position.
$ws.Range("D3").Value = 34

# Best-effort: also nudge the control's own cached Value/position through
# the Forms-control object model so any engine that does track it stays in
# sync with the linked cell (no-op on engines that only track this via the
# linked cell itself).
foreach ($shp in $ws.Shapes) {
    if ($shp.Name -eq "Scroll Bar 1") {
        $cf = $shp.ControlFormat
        $cf.Value = 34
    }
}

$wb.Application.Calculate()
